$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (rich text shared strings) ---
$ws.Range("A8").Value = "Volume 30   Number  51"
$ws.Range("C9").Value = "Report Covering the Week  12/18/2023  Through  12/24/2023"

# --- Data table updates ---
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$ws.Range("C15").Copy() | Out-Null
$ws.Range("D15").PasteSpecial(-4122) | Out-Null
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$ws.Range("F15").Copy() | Out-Null
$ws.Range("E15").PasteSpecial(-4122) | Out-Null

# Row 16
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 66.666666666666
$ws.Range("F16").Value = 16
$ws.Range("G16").Value = 8
$ws.Range("I16").Value = 100
$ws.Range("J16").Value = 115
$ws.Range("K16").Value = -13.043478260869
$ws.Range("L16").Value = -9.090909090909
$ws.Range("M16").Value = 12.359550561797
$ws.Range("N16").Value = -83.108108108108

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0"
$ws.Range("A17").Copy() | Out-Null
$ws.Range("D17").PasteSpecial(-4122) | Out-Null
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "***.*"
$ws.Range("A17").Copy() | Out-Null
$ws.Range("E17").PasteSpecial(-4122) | Out-Null
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 4
$ws.Range("H17").Value = 150
$ws.Range("I17").Value = 107
$ws.Range("K17").Value = 7
$ws.Range("L17").Value = 24.418604651162
$ws.Range("M17").Value = 91.071428571428
$ws.Range("N17").Value = 13.829787234042

# Row 18
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 1
$ws.Range("E18").Value = 100
$ws.Range("I18").Value = 107
$ws.Range("J18").Value = 135
$ws.Range("K18").Value = -20.74074074074
$ws.Range("L18").Value = 35.443037974683
$ws.Range("M18").Value = 7
$ws.Range("N18").Value = -88.0978865406

# Row 19
$ws.Range("C19").Value = 16
$ws.Range("D19").Value = 16
$ws.Range("F19").Value = 69
$ws.Range("G19").Value = 51
$ws.Range("H19").Value = 35.294117647058
$ws.Range("I19").Value = 737
$ws.Range("J19").Value = 804
$ws.Range("K19").Value = -8.333333333333
$ws.Range("L19").Value = 1.375515818431
$ws.Range("M19").Value = 9.347181008902
$ws.Range("N19").Value = -58.43203609701

# Row 20
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 2
$ws.Range("I20").Value = 92
$ws.Range("K20").Value = 31.428571428571
$ws.Range("L20").Value = 6.976744186046
$ws.Range("M20").Value = 178.787878787879
$ws.Range("N20").Value = -91.908531222515

# Row 21
$ws.Range("C21").Value = 27
$ws.Range("D21").Value = 20
$ws.Range("E21").Value = 35
$ws.Range("F21").Value = 107
$ws.Range("G21").Value = 78
$ws.Range("H21").Value = 37.179487179487
$ws.Range("I21").Value = 1157
$ws.Range("J21").Value = 1238
$ws.Range("K21").Value = -6.54281098546
$ws.Range("L21").Value = 5.662100456621
$ws.Range("M21").Value = 20.020746887966
$ws.Range("N21").Value = -74.379982285208

# Row 22
$ws.Range("F22").Value = 4
$ws.Range("G22").Value = 4
$ws.Range("I22").Value = 28
$ws.Range("J22").Value = 27
$ws.Range("K22").Value = 3.703703703703
$ws.Range("L22").Value = 21.739130434782
$ws.Range("M22").Value = -15.151515151515

# Row 23
$ws.Range("C23").Value = 3
$ws.Range("F23").Value = 8
$ws.Range("I23").Value = 43
$ws.Range("K23").Value = 79.166666666666
$ws.Range("L23").Value = 72
$ws.Range("M23").Value = 65.384615384615

# Row 24
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = -8.333333333333
$ws.Range("F24").Value = 96
$ws.Range("G24").Value = 70
$ws.Range("H24").Value = 37.142857142857
$ws.Range("I24").Value = 1158
$ws.Range("J24").Value = 1205
$ws.Range("K24").Value = -3.900414937759
$ws.Range("L24").Value = -6.006493506493
$ws.Range("M24").Value = 12.536443148688

# Row 25
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = 33.333333333333
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 17
$ws.Range("H25").Value = -11.764705882352
$ws.Range("I25").Value = 220
$ws.Range("J25").Value = 222
$ws.Range("K25").Value = -0.9009009009
$ws.Range("L25").Value = 16.402116402116
$ws.Range("M25").Value = -19.70802919708

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0"
$ws.Range("C26").Copy() | Out-Null
$ws.Range("D26").PasteSpecial(-4122) | Out-Null
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "***.*"
$ws.Range("F26").Copy() | Out-Null
$ws.Range("E26").PasteSpecial(-4122) | Out-Null

# Row 27
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$ws.Range("M27").Copy() | Out-Null
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$ws.Range("M27").Copy() | Out-Null
$ws.Range("D27").PasteSpecial(-4122) | Out-Null
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$ws.Range("M27").Copy() | Out-Null
$ws.Range("E27").PasteSpecial(-4122) | Out-Null
$ws.Range("F27").Value = 2
$ws.Range("H27").Value = -60
$ws.Range("L27").Value = -24.074074074074

# Row 28
$ws.Range("C28").Value = 1
$ws.Range("I28").Copy() | Out-Null
$ws.Range("C28").PasteSpecial(-4122) | Out-Null
$ws.Range("C28").Value = 1
$ws.Range("F28").Value = 1
$ws.Range("I28").Copy() | Out-Null
$ws.Range("F28").PasteSpecial(-4122) | Out-Null
$ws.Range("F28").Value = 1
$ws.Range("I28").Value = 3
$ws.Range("K28").Value = 200
$ws.Range("L28").Value = 0
$ws.Range("N28").Value = 50

# Row 29
$ws.Range("C29").Value = 1
$ws.Range("I29").Copy() | Out-Null
$ws.Range("C29").PasteSpecial(-4122) | Out-Null
$ws.Range("C29").Value = 1
$ws.Range("F29").Value = 1
$ws.Range("I29").Copy() | Out-Null
$ws.Range("F29").PasteSpecial(-4122) | Out-Null
$ws.Range("F29").Value = 1
$ws.Range("I29").Value = 3
$ws.Range("K29").Value = 200
$ws.Range("L29").Value = 0
$ws.Range("N29").Value = 50

$excel.CutCopyMode = $false